$wb = $excel.ActiveWorkbook

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3719.0938
$ws.Range("I64").Value = 3637.375
$ws.Range("J64").Value = 3800.8125
$ws.Range("K64").Value = 3637.375
$ws.Range("L64").Value = 3800.8125
$ws.Range("M64").Value = -3389.375
$ws.Range("N64").Value = -4296.8125

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3719.0938
$ws.Range("I67").Value = 3637.375
$ws.Range("J67").Value = 3800.8125
$ws.Range("K67").Value = 3637.375
$ws.Range("L67").Value = 3800.8125
$ws.Range("M67").Value = -2779.375
$ws.Range("N67").Value = -5516.8125

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 7707.4517
$ws.Range("I76").Value = 9958.733
$ws.Range("J76").Value = 5596.875
$ws.Range("K76").Value = 9958.733
$ws.Range("L76").Value = 5596.875
$ws.Range("M76").Value = -9643.733
$ws.Range("N76").Value = -6226.875

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 7707.4517
$ws.Range("I79").Value = 9958.733
$ws.Range("J79").Value = 5596.875
$ws.Range("K79").Value = 9958.733
$ws.Range("L79").Value = 5596.875
$ws.Range("M79").Value = -8866.733
$ws.Range("N79").Value = -7780.875

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 372.27274
$ws.Range("I107").Value = 366.1111
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 366.1111
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1553.8889
$ws.Range("N107").Value = -4240

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1418.25
$ws.Range("I45").Value = 1177.5
$ws.Range("J45").Value = 1899.75
$ws.Range("K45").Value = 1177.5
$ws.Range("L45").Value = 1899.75
$ws.Range("M45").Value = -800.5
$ws.Range("N45").Value = -2653.75

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2545.5103
$ws.Range("I105").Value = 2363
$ws.Range("J105").Value = 3356.6667
$ws.Range("K105").Value = 2363
$ws.Range("L105").Value = 3356.6667
$ws.Range("M105").Value = -616
$ws.Range("N105").Value = -6850.6667

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 4653
$ws.Range("I41").Value = 4653
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 4653
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -4225
$ws.Range("N41").ClearContents()

# CRP row 50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

# CRP row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 24756
$ws.Range("J51").Value = 24756
$ws.Range("L51").Value = 24756
$ws.Range("N51").Value = -26228

# CRP row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 24725
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 29633.334
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 29633.334
$ws.Range("M60").Value = -9489
$ws.Range("N60").Value = -30655.334

# CRP row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 24756
$ws.Range("J61").Value = 24756
$ws.Range("L61").Value = 24756
$ws.Range("N61").Value = -25452

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3730.4614
$ws.Range("I62").Value = 3506.0605
$ws.Range("J62").Value = 4964.6665
$ws.Range("K62").Value = 3506.0605
$ws.Range("L62").Value = 4964.6665
$ws.Range("M62").Value = -2882.0605
$ws.Range("N62").Value = -6212.6665

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3730.4614
$ws.Range("I65").Value = 3506.0605
$ws.Range("J65").Value = 4964.6665
$ws.Range("K65").Value = 17530.3025
$ws.Range("L65").Value = 24823.3325
$ws.Range("M65").Value = -14410.3025
$ws.Range("N65").Value = -31063.3325

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 451.02563
$ws.Range("I5").Value = 311.5
$ws.Range("K5").Value = 934.5
$ws.Range("M5").Value = -822.5

# CUL row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 608.8889
$ws.Range("I7").Value = 82.5
$ws.Range("J7").Value = 759.2857
$ws.Range("K7").Value = 247.5
$ws.Range("L7").Value = 2277.8571
$ws.Range("M7").Value = -135.5
$ws.Range("N7").Value = -2501.8571

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 451.02563
$ws.Range("I135").Value = 311.5
$ws.Range("K135").Value = 2803.5
$ws.Range("M135").Value = -268.5

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7269.811
$ws.Range("I70").Value = 3864.9656
$ws.Range("J70").Value = 19612.375
$ws.Range("K70").Value = 3864.9656
$ws.Range("L70").Value = 19612.375
$ws.Range("M70").Value = -3594.9656
$ws.Range("N70").Value = -20152.375

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7269.811
$ws.Range("I73").Value = 3864.9656
$ws.Range("J73").Value = 19612.375
$ws.Range("K73").Value = 3864.9656
$ws.Range("L73").Value = 19612.375
$ws.Range("M73").Value = -2928.9656
$ws.Range("N73").Value = -21484.375

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4031.8572
$ws.Range("I80").Value = 4584.1304
$ws.Range("J80").Value = 2973.3333
$ws.Range("K80").Value = 4584.1304
$ws.Range("L80").Value = 2973.3333
$ws.Range("M80").Value = -3586.1304
$ws.Range("N80").Value = -4969.3333

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4031.8572
$ws.Range("I83").Value = 4584.1304
$ws.Range("J83").Value = 2973.3333
$ws.Range("K83").Value = 22920.652
$ws.Range("L83").Value = 14866.6665
$ws.Range("M83").Value = -17928.652
$ws.Range("N83").Value = -24850.6665

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1113.0588
$ws.Range("I122").Value = 1008.7143
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 3026.1429
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -576.1428999999998
$ws.Range("N122").Value = -9700

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2589.7693
$ws.Range("I7").Value = 2525
$ws.Range("J7").Value = 2665.3333
$ws.Range("K7").Value = 2525
$ws.Range("L7").Value = 2665.3333
$ws.Range("M7").Value = -2413
$ws.Range("N7").Value = -2889.3333

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2374.8948
$ws.Range("I40").Value = 2456.6365
$ws.Range("J40").Value = 2262.5
$ws.Range("K40").Value = 2456.6365
$ws.Range("L40").Value = 2262.5
$ws.Range("M40").Value = -2320.6365
$ws.Range("N40").Value = -2534.5

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2589.7693
$ws.Range("I126").Value = 2525
$ws.Range("J126").Value = 2665.3333
$ws.Range("K126").Value = 7575
$ws.Range("L126").Value = 7995.999899999999
$ws.Range("M126").Value = -5105
$ws.Range("N126").Value = -12935.9999

# WVR row 15
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 4183.3335
$ws.Range("J15").Value = 4183.3335
$ws.Range("L15").Value = 4183.3335
$ws.Range("N15").Value = -4759.3335

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 35179.566
$ws.Range("I122").Value = 41707.28
$ws.Range("J122").Value = 2541
$ws.Range("K122").Value = 125121.84
$ws.Range("L122").Value = 7623
$ws.Range("M122").Value = -122671.84
$ws.Range("N122").Value = -12523

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2498.8333
$ws.Range("I132").Value = 1430.7333
$ws.Range("J132").Value = 3092.2222
$ws.Range("K132").Value = 4292.199900000001
$ws.Range("L132").Value = 9276.6666
$ws.Range("M132").Value = -1762.199900000001
$ws.Range("N132").Value = -14336.6666
